$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

# Replace $old with $new inside a single table cell. Cell.Range itself is
# not reliably bounded by this host's Find engine (it behaves like a
# document-wide search), so a fresh Range rebuilt from the cell's
# Start/End offsets plus Wrap:=wdFindStop (0) and Replace:=wdReplaceOne (1)
# is used to keep the match (and the edit) inside the target cell only.
function ReplaceInCell($table, $row, $col, $old, $new) {
    $cell = $table.Cell($row, $col)
    $cellRange = $cell.Range
    $scoped = $d.Range($cellRange.Start, $cellRange.End)
    $scoped.Find.Execute($old, $false, $false, $false, $false, $false, $true, 0, $false, $new, 1) | Out-Null
}

# Row 2 ("dateTime"): only the "Champ correspondant" column changes.
ReplaceInCell $t 2 2 "Date Heure de création de la demande" "Date Heure de réponse"

# Row 3 ("convention" -> "answer"): tag name, field, format (+ enum line),
# cardinality and description all change.
ReplaceInCell $t 3 1 "convention" "answer"
ReplaceInCell $t 3 2 "Cadre conventionnel" "Réponse"
ReplaceInCell $t 3 3 "string" ("string" + [char]11 + "(ENUM : OUI, NON, PARTIEL, DIFFERE)")
ReplaceInCell $t 3 4 "0..1" "1..1"
ReplaceInCell $t 3 5 "Nomenclature ? " "oui / non / oui partiel / différé"

# Row 4 ("deadline"): field and description change.
ReplaceInCell $t 4 2 "Délai souhaité" "Délai de réponse"
# This description contains apostrophes; Find/Replace would smart-quote
# them, so set the cell text directly instead to keep plain "'" marks.
$t.Cell(4, 5).Range.Text = "Indique le délai de réponse auquel s'engage l'expéditeur"

# Row 5 ("purpose") is removed entirely.
$t.Rows.Item(5).Delete()

# Former row 6 ("freetext") is now row 5: field and description change.
ReplaceInCell $t 5 2 "Précisions sur la demande" "Précisions sur la réponse"
ReplaceInCell $t 5 5 "Texte libre permettant de détailler la demande" "Commentaire libre pour apporter toutes précisions utiles à la réponse"
